$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ptgs2"
$ws.Range("C2").Value = "Cav1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8976106666666667
$ws.Range("H2").Value = 2.692832
$ws.Range("I2").Value = 0.02124022359297476
$ws.Range("J2").Value = 0.02124022359297476
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 833.4308676666666
$ws.Range("N2").Value = 2500.292603
$ws.Range("O2").Value = 0.8518935545813505
$ws.Range("P2").Value = 0.8518935545813505
$ws.Range("Q2").Value = 748.096436746855
$ws.Range("R2").Value = 6732.867930721696
$ws.Range("S2").Value = 0.01809440957672193
$ws.Range("T2").Value = 0.01809440957672193

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ptgs2"
$ws.Range("C3").Value = "Cav1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8976106666666667
$ws.Range("H3").Value = 2.692832
$ws.Range("I3").Value = 0.02124022359297476
$ws.Range("J3").Value = 0.02124022359297476
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.73945766666667
$ws.Range("N3").Value = 53.218373
$ws.Range("O3").Value = 0.01813243333584592
$ws.Range("P3").Value = 0.01813243333584592
$ws.Range("Q3").Value = 15.92312642248178
$ws.Range("R3").Value = 143.308137802336
$ws.Range("S3").Value = 0.0003851369383380766
$ws.Range("T3").Value = 0.0003851369383380766

# Row 4: ECs -> M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ptgs2"
$ws.Range("C4").Value = "Cav1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8976106666666667
$ws.Range("H4").Value = 2.692832
$ws.Range("I4").Value = 0.02124022359297476
$ws.Range("J4").Value = 0.02124022359297476
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.344749666666667
$ws.Range("N4").Value = 4.034249
$ws.Range("O4").Value = 0.001374539410528448
$ws.Range("P4").Value = 0.001374539410528448
$ws.Range("Q4").Value = 1.207061644796444
$ws.Range("R4").Value = 10.863554803168
$ws.Range("S4").Value = 0.00002919552441697997
$ws.Range("T4").Value = 0.00002919552441697996

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ptgs2"
$ws.Range("C5").Value = "Cav1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8976106666666667
$ws.Range("H5").Value = 2.692832
$ws.Range("I5").Value = 0.02124022359297476
$ws.Range("J5").Value = 0.02124022359297476
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 125.812397
$ws.Range("N5").Value = 377.437191
$ws.Range("O5").Value = 0.1285994726722751
$ws.Range("P5").Value = 0.1285994726722751
$ws.Range("Q5").Value = 112.9305495461013
$ws.Range("R5").Value = 1016.374945914912
$ws.Range("S5").Value = 0.002731481553497772
$ws.Range("T5").Value = 0.002731481553497772

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ptgs2"
$ws.Range("C6").Value = "Cav1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 20.60823266666667
$ws.Range("H6").Value = 61.824698
$ws.Range("I6").Value = 0.4876540419484541
$ws.Range("J6").Value = 0.4876540419484541
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 833.4308676666666
$ws.Range("N6").Value = 2500.292603
$ws.Range("O6").Value = 0.8518935545813505
$ws.Range("P6").Value = 0.8518935545813505
$ws.Range("Q6").Value = 17175.53723245654
$ws.Range("R6").Value = 154579.8350921089
$ws.Range("S6").Value = 0.4154293352014315
$ws.Range("T6").Value = 0.4154293352014315

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ptgs2"
$ws.Range("C7").Value = "Cav1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 20.60823266666667
$ws.Range("H7").Value = 61.824698
$ws.Range("I7").Value = 0.4876540419484541
$ws.Range("J7").Value = 0.4876540419484541
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 17.73945766666667
$ws.Range("N7").Value = 53.218373
$ws.Range("O7").Value = 0.01813243333584592
$ws.Range("P7").Value = 0.01813243333584592
$ws.Range("Q7").Value = 365.5788709751504
$ws.Range("R7").Value = 3290.209838776354
$ws.Range("S7").Value = 0.008842354406586154
$ws.Range("T7").Value = 0.008842354406586154

# Row 8: FAPs -> M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Ptgs2"
$ws.Range("C8").Value = "Cav1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 20.60823266666667
$ws.Range("H8").Value = 61.824698
$ws.Range("I8").Value = 0.4876540419484541
$ws.Range("J8").Value = 0.4876540419484541
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.344749666666667
$ws.Range("N8").Value = 4.034249
$ws.Range("O8").Value = 0.001374539410528448
$ws.Range("P8").Value = 0.001374539410528448
$ws.Range("Q8").Value = 27.71291400908911
$ws.Range("R8").Value = 249.416226081802
$ws.Range("S8").Value = 0.0006702996993616432
$ws.Range("T8").Value = 0.0006702996993616431

# Row 9: FAPs -> sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Ptgs2"
$ws.Range("C9").Value = "Cav1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 20.60823266666667
$ws.Range("H9").Value = 61.824698
$ws.Range("I9").Value = 0.4876540419484541
$ws.Range("J9").Value = 0.4876540419484541
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 125.812397
$ws.Range("N9").Value = 377.437191
$ws.Range("O9").Value = 0.1285994726722751
$ws.Range("P9").Value = 0.1285994726722751
$ws.Range("Q9").Value = 2592.771149727036
$ws.Range("R9").Value = 23334.94034754332
$ws.Range("S9").Value = 0.06271205264107474
$ws.Range("T9").Value = 0.06271205264107474

# Row 10: M2 -> ECs
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Ptgs2"
$ws.Range("C10").Value = "Cav1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 20.573189
$ws.Range("H10").Value = 61.719567
$ws.Range("I10").Value = 0.4868248012284415
$ws.Range("J10").Value = 0.4868248012284414
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 833.4308676666666
$ws.Range("N10").Value = 2500.292603
$ws.Range("O10").Value = 0.8518935545813505
$ws.Range("P10").Value = 0.8518935545813505
$ws.Range("Q10").Value = 17146.33075894032
$ws.Range("R10").Value = 154316.9768304629
$ws.Range("S10").Value = 0.4147229103768564
$ws.Range("T10").Value = 0.4147229103768564

# Row 11: M2 -> FAPs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Ptgs2"
$ws.Range("C11").Value = "Cav1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 20.573189
$ws.Range("H11").Value = 61.719567
$ws.Range("I11").Value = 0.4868248012284415
$ws.Range("J11").Value = 0.4868248012284414
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 17.73945766666667
$ws.Range("N11").Value = 53.218373
$ws.Range("O11").Value = 0.01813243333584592
$ws.Range("P11").Value = 0.01813243333584592
$ws.Range("Q11").Value = 364.9572153338323
$ws.Range("R11").Value = 3284.614938004491
$ws.Range("S11").Value = 0.008827318254511156
$ws.Range("T11").Value = 0.008827318254511154

# Row 12: M2 -> M2
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Ptgs2"
$ws.Range("C12").Value = "Cav1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 20.573189
$ws.Range("H12").Value = 61.719567
$ws.Range("I12").Value = 0.4868248012284415
$ws.Range("J12").Value = 0.4868248012284414
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.344749666666667
$ws.Range("N12").Value = 4.034249
$ws.Range("O12").Value = 0.001374539410528448
$ws.Range("P12").Value = 0.001374539410528448
$ws.Range("Q12").Value = 27.66578905002033
$ws.Range("R12").Value = 248.992101450183
$ws.Range("S12").Value = 0.0006691598753111709
$ws.Range("T12").Value = 0.0006691598753111707

# Row 13: M2 -> sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Ptgs2"
$ws.Range("C13").Value = "Cav1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 20.573189
$ws.Range("H13").Value = 61.719567
$ws.Range("I13").Value = 0.4868248012284415
$ws.Range("J13").Value = 0.4868248012284414
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 125.812397
$ws.Range("N13").Value = 377.437191
$ws.Range("O13").Value = 0.1285994726722751
$ws.Range("P13").Value = 0.1285994726722751
$ws.Range("Q13").Value = 2588.362222024033
$ws.Range("R13").Value = 23295.2599982163
$ws.Range("S13").Value = 0.06260541272176273
$ws.Range("T13").Value = 0.06260541272176273

# Row 14: sCs -> ECs
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Ptgs2"
$ws.Range("C14").Value = "Cav1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.180912
$ws.Range("H14").Value = 0.542736
$ws.Range("I14").Value = 0.004280933230129748
$ws.Range("J14").Value = 0.004280933230129748
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 833.4308676666666
$ws.Range("N14").Value = 2500.292603
$ws.Range("O14").Value = 0.8518935545813505
$ws.Range("P14").Value = 0.8518935545813505
$ws.Range("Q14").Value = 150.777645131312
$ws.Range("R14").Value = 1356.998806181808
$ws.Range("S14").Value = 0.003646899426340654
$ws.Range("T14").Value = 0.003646899426340653

# Row 15: sCs -> FAPs
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Ptgs2"
$ws.Range("C15").Value = "Cav1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.180912
$ws.Range("H15").Value = 0.542736
$ws.Range("I15").Value = 0.004280933230129748
$ws.Range("J15").Value = 0.004280933230129748
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 17.73945766666667
$ws.Range("N15").Value = 53.218373
$ws.Range("O15").Value = 0.01813243333584592
$ws.Range("P15").Value = 0.01813243333584592
$ws.Range("Q15").Value = 3.209280765392
$ws.Range("R15").Value = 28.883526888528
$ws.Range("S15").Value = 0.00007762373641053521
$ws.Range("T15").Value = 0.00007762373641053518

# Row 16: sCs -> M2
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Ptgs2"
$ws.Range("C16").Value = "Cav1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.180912
$ws.Range("H16").Value = 0.542736
$ws.Range("I16").Value = 0.004280933230129748
$ws.Range("J16").Value = 0.004280933230129748
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.344749666666667
$ws.Range("N16").Value = 4.034249
$ws.Range("O16").Value = 0.001374539410528448
$ws.Range("P16").Value = 0.001374539410528448
$ws.Range("Q16").Value = 0.243281351696
$ws.Range("R16").Value = 2.189532165264
$ws.Range("S16").Value = 0.00000588431143865419
$ws.Range("T16").Value = 0.000005884311438654188

# Row 17: sCs -> sCs
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Ptgs2"
$ws.Range("C17").Value = "Cav1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.180912
$ws.Range("H17").Value = 0.542736
$ws.Range("I17").Value = 0.004280933230129748
$ws.Range("J17").Value = 0.004280933230129748
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 125.812397
$ws.Range("N17").Value = 377.437191
$ws.Range("O17").Value = 0.1285994726722751
$ws.Range("P17").Value = 0.1285994726722751
$ws.Range("Q17").Value = 22.760972366064
$ws.Range("R17").Value = 204.848751294576
$ws.Range("S17").Value = 0.0005505257559399051
$ws.Range("T17").Value = 0.000550525755939905
